$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Pin" column (B) values for rows 2-7
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = 7

# Fill in the previously-empty Min (C) and Max (E) cells for row 7
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 180

# Update the active cell selection on the sheet
[void]$ws.Range("E8").Select()
